$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24
$ws.Cells.Item(24, 1).Value = 112390031
$ws.Cells.Item(24, 2).Value = 90826
$ws.Cells.Item(24, 4).Value = 'LC'
$ws.Cells.Item(24, 5).Value = 4366
$ws.Cells.Item(24, 6).Value = 'Skarp dropptaggsvamp'
$ws.Cells.Item(24, 7).Value = 'Hydnellum peckii'
$ws.Cells.Item(24, 8).Value = 'Banker'
$ws.Cells.Item(24, 9).NumberFormat = "@"
$ws.Cells.Item(24, 9).Value = '1'
$ws.Cells.Item(24, 10).ClearContents()
$ws.Cells.Item(24, 17).Value = 509098
$ws.Cells.Item(24, 18).Value = 6784229
$ws.Cells.Item(24, 34).Value = 'Sandtallskog'

# Row 25
$ws.Cells.Item(25, 1).Value = 112390262
$ws.Cells.Item(25, 2).Value = 90808
$ws.Cells.Item(25, 4).Value = 'NT'
$ws.Cells.Item(25, 5).Value = 4362
$ws.Cells.Item(25, 6).Value = 'Blå taggsvamp'
$ws.Cells.Item(25, 7).Value = 'Hydnellum caeruleum'
$ws.Cells.Item(25, 8).Value = '(Hornem.) P.Karst.'
$ws.Cells.Item(25, 9).ClearContents()
$ws.Cells.Item(25, 10).ClearContents()
$ws.Cells.Item(25, 17).Value = 509072
$ws.Cells.Item(25, 18).Value = 6784116
$ws.Cells.Item(25, 34).Value = 'Sandtallskog'

# Row 26
$ws.Cells.Item(26, 1).Value = 112390292
$ws.Cells.Item(26, 2).Value = 90808
$ws.Cells.Item(26, 4).Value = 'NT'
$ws.Cells.Item(26, 5).Value = 4362
$ws.Cells.Item(26, 6).Value = 'Blå taggsvamp'
$ws.Cells.Item(26, 7).Value = 'Hydnellum caeruleum'
$ws.Cells.Item(26, 8).Value = '(Hornem.) P.Karst.'
$ws.Cells.Item(26, 9).NumberFormat = "@"
$ws.Cells.Item(26, 9).Value = '1'
$ws.Cells.Item(26, 10).Value = 'fruktkroppar'
$ws.Cells.Item(26, 17).Value = 509065
$ws.Cells.Item(26, 18).Value = 6784066
$ws.Cells.Item(26, 34).ClearContents()

# Row 27
$ws.Cells.Item(27, 1).Value = 112390287
$ws.Cells.Item(27, 2).Value = 89072
$ws.Cells.Item(27, 4).Value = 'LC'
$ws.Cells.Item(27, 5).Value = 256703
$ws.Cells.Item(27, 6).Value = 'Tallfingersvamp'
$ws.Cells.Item(27, 7).Value = 'Ramaria eosanguinea'
$ws.Cells.Item(27, 8).Value = 'R.H.Petersen'
$ws.Cells.Item(27, 9).ClearContents()
$ws.Cells.Item(27, 10).ClearContents()
$ws.Cells.Item(27, 17).Value = 509070
$ws.Cells.Item(27, 18).Value = 6784097
$ws.Cells.Item(27, 34).Value = 'Sandtallskog'

# Row 28
$ws.Cells.Item(28, 1).Value = 112390426
$ws.Cells.Item(28, 2).Value = 90826
$ws.Cells.Item(28, 4).Value = 'LC'
$ws.Cells.Item(28, 5).Value = 4366
$ws.Cells.Item(28, 6).Value = 'Skarp dropptaggsvamp'
$ws.Cells.Item(28, 7).Value = 'Hydnellum peckii'
$ws.Cells.Item(28, 8).Value = 'Banker'
$ws.Cells.Item(28, 9).NumberFormat = "@"
$ws.Cells.Item(28, 9).Value = '1'
$ws.Cells.Item(28, 10).Value = 'fruktkroppar'
$ws.Cells.Item(28, 17).Value = 509076
$ws.Cells.Item(28, 18).Value = 6783959
$ws.Cells.Item(28, 34).Value = 'Sandtallskog'

# Row 29
$ws.Cells.Item(29, 1).Value = 112390524
$ws.Cells.Item(29, 2).Value = 90814
$ws.Cells.Item(29, 4).Value = 'LC'
$ws.Cells.Item(29, 5).Value = 4364
$ws.Cells.Item(29, 6).Value = 'Dropptaggsvamp'
$ws.Cells.Item(29, 7).Value = 'Hydnellum ferrugineum'
$ws.Cells.Item(29, 8).Value = '(Fr.:Fr.) P. Karst.'
$ws.Cells.Item(29, 9).ClearContents()
$ws.Cells.Item(29, 10).ClearContents()
$ws.Cells.Item(29, 17).Value = 509060
$ws.Cells.Item(29, 18).Value = 6783866
$ws.Cells.Item(29, 34).ClearContents()

# Row 30
$ws.Cells.Item(30, 1).Value = 112389959
$ws.Cells.Item(30, 2).Value = 89072
$ws.Cells.Item(30, 4).Value = 'LC'
$ws.Cells.Item(30, 5).Value = 256703
$ws.Cells.Item(30, 6).Value = 'Tallfingersvamp'
$ws.Cells.Item(30, 7).Value = 'Ramaria eosanguinea'
$ws.Cells.Item(30, 8).Value = 'R.H.Petersen'
$ws.Cells.Item(30, 9).NumberFormat = "@"
$ws.Cells.Item(30, 9).Value = '1'
$ws.Cells.Item(30, 10).Value = 'fruktkroppar'
$ws.Cells.Item(30, 17).Value = 509111
$ws.Cells.Item(30, 18).Value = 6784257
$ws.Cells.Item(30, 34).Value = 'Sandtallskog'

# Row 31
$ws.Cells.Item(31, 1).Value = 112389988
$ws.Cells.Item(31, 2).Value = 89072
$ws.Cells.Item(31, 4).Value = 'LC'
$ws.Cells.Item(31, 5).Value = 256703
$ws.Cells.Item(31, 6).Value = 'Tallfingersvamp'
$ws.Cells.Item(31, 7).Value = 'Ramaria eosanguinea'
$ws.Cells.Item(31, 8).Value = 'R.H.Petersen'
$ws.Cells.Item(31, 9).NumberFormat = "@"
$ws.Cells.Item(31, 9).Value = '1'
$ws.Cells.Item(31, 10).ClearContents()
$ws.Cells.Item(31, 17).Value = 509101
$ws.Cells.Item(31, 18).Value = 6784234
$ws.Cells.Item(31, 34).Value = 'Sandtallskog'

# Row 32
$ws.Cells.Item(32, 1).Value = 112390256
$ws.Cells.Item(32, 2).Value = 90448
$ws.Cells.Item(32, 4).Value = 'NT'
$ws.Cells.Item(32, 5).Value = 4745
$ws.Cells.Item(32, 6).Value = 'Tallriska'
$ws.Cells.Item(32, 7).Value = 'Lactarius musteus'
$ws.Cells.Item(32, 8).Value = 'Fr.'
$ws.Cells.Item(32, 9).NumberFormat = "@"
$ws.Cells.Item(32, 9).Value = '1'
$ws.Cells.Item(32, 10).Value = 'fruktkroppar'
$ws.Cells.Item(32, 17).Value = 509090
$ws.Cells.Item(32, 18).Value = 6784191
$ws.Cells.Item(32, 34).Value = 'Sandtallskog'

# Row 33
$ws.Cells.Item(33, 1).Value = 112390451
$ws.Cells.Item(33, 2).Value = 90814
$ws.Cells.Item(33, 4).Value = 'LC'
$ws.Cells.Item(33, 5).Value = 4364
$ws.Cells.Item(33, 6).Value = 'Dropptaggsvamp'
$ws.Cells.Item(33, 7).Value = 'Hydnellum ferrugineum'
$ws.Cells.Item(33, 8).Value = '(Fr.:Fr.) P. Karst.'
$ws.Cells.Item(33, 9).NumberFormat = "@"
$ws.Cells.Item(33, 9).Value = '3'
$ws.Cells.Item(33, 10).Value = 'fruktkroppar'
$ws.Cells.Item(33, 17).Value = 509076
$ws.Cells.Item(33, 18).Value = 6783959
$ws.Cells.Item(33, 34).Value = 'Sandtallskog'

# Row 34
$ws.Cells.Item(34, 1).Value = 112390567
$ws.Cells.Item(34, 2).Value = 90806
$ws.Cells.Item(34, 4).Value = 'NT'
$ws.Cells.Item(34, 5).Value = 4361
$ws.Cells.Item(34, 6).Value = 'Orange taggsvamp'
$ws.Cells.Item(34, 7).Value = 'Hydnellum aurantiacum'
$ws.Cells.Item(34, 8).Value = '(Batsch:Fr.) P.Karst.'
$ws.Cells.Item(34, 9).NumberFormat = "@"
$ws.Cells.Item(34, 9).Value = '19'
$ws.Cells.Item(34, 10).Value = 'fruktkroppar'
$ws.Cells.Item(34, 17).Value = 509010
$ws.Cells.Item(34, 18).Value = 6783836
$ws.Cells.Item(34, 34).Value = 'Sandtallskog'

# Row 35
$ws.Cells.Item(35, 1).Value = 112390119
$ws.Cells.Item(35, 2).Value = 90830
$ws.Cells.Item(35, 4).Value = 'NT'
$ws.Cells.Item(35, 5).Value = 2059
$ws.Cells.Item(35, 6).Value = 'Skrovlig taggsvamp'
$ws.Cells.Item(35, 7).Value = 'Hydnellum scabrosum'
$ws.Cells.Item(35, 8).Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Cells.Item(35, 9).ClearContents()
$ws.Cells.Item(35, 10).ClearContents()
$ws.Cells.Item(35, 17).Value = 509093
$ws.Cells.Item(35, 18).Value = 6784215
$ws.Cells.Item(35, 34).Value = 'Sandtallskog'

# Row 36
$ws.Cells.Item(36, 1).Value = 112390398
$ws.Cells.Item(36, 2).Value = 90830
$ws.Cells.Item(36, 4).Value = 'NT'
$ws.Cells.Item(36, 5).Value = 2059
$ws.Cells.Item(36, 6).Value = 'Skrovlig taggsvamp'
$ws.Cells.Item(36, 7).Value = 'Hydnellum scabrosum'
$ws.Cells.Item(36, 8).Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Cells.Item(36, 9).ClearContents()
$ws.Cells.Item(36, 10).ClearContents()
$ws.Cells.Item(36, 17).Value = 509066
$ws.Cells.Item(36, 18).Value = 6784010
$ws.Cells.Item(36, 34).Value = 'Sandtallskog'

# Row 37
$ws.Cells.Item(37, 1).Value = 112390630
$ws.Cells.Item(37, 2).Value = 90857
$ws.Cells.Item(37, 4).Value = 'NT'
$ws.Cells.Item(37, 5).Value = 5448
$ws.Cells.Item(37, 6).Value = 'Svartvit taggsvamp'
$ws.Cells.Item(37, 7).Value = 'Phellodon connatus'
$ws.Cells.Item(37, 8).Value = '(Schultz) nom.prov'
$ws.Cells.Item(37, 9).ClearContents()
$ws.Cells.Item(37, 10).ClearContents()
$ws.Cells.Item(37, 17).Value = 509014
$ws.Cells.Item(37, 18).Value = 6783848
$ws.Cells.Item(37, 34).ClearContents()

# Row 38
$ws.Cells.Item(38, 1).Value = 112390509
$ws.Cells.Item(38, 2).Value = 90830
$ws.Cells.Item(38, 4).Value = 'NT'
$ws.Cells.Item(38, 5).Value = 2059
$ws.Cells.Item(38, 6).Value = 'Skrovlig taggsvamp'
$ws.Cells.Item(38, 7).Value = 'Hydnellum scabrosum'
$ws.Cells.Item(38, 8).Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Cells.Item(38, 9).ClearContents()
$ws.Cells.Item(38, 10).ClearContents()
$ws.Cells.Item(38, 17).Value = 509056
$ws.Cells.Item(38, 18).Value = 6783885
$ws.Cells.Item(38, 34).ClearContents()

# Row 39
$ws.Cells.Item(39, 1).Value = 112390382
$ws.Cells.Item(39, 2).Value = 90830
$ws.Cells.Item(39, 4).Value = 'NT'
$ws.Cells.Item(39, 5).Value = 2059
$ws.Cells.Item(39, 6).Value = 'Skrovlig taggsvamp'
$ws.Cells.Item(39, 7).Value = 'Hydnellum scabrosum'
$ws.Cells.Item(39, 8).Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Cells.Item(39, 9).ClearContents()
$ws.Cells.Item(39, 10).ClearContents()
$ws.Cells.Item(39, 17).Value = 509061
$ws.Cells.Item(39, 18).Value = 6784061
$ws.Cells.Item(39, 34).Value = 'Sandtallskog'
